$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.220.17"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.577.99"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'571.54"
$ws.Range("E5").Value = "  +2.86%  "
$ws.Range("D6").Value = "'143.19"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "2.583.87"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "'6.70"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").Value = "'0.103"
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("E12").Value = "  +11.78%  "
$ws.Range("D13").Value = "'0.345"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").Value = "3.028.19"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "59.235.49"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "'22.51"
$ws.Range("E16").Value = "  +7.67%  "
$ws.Range("E17").Value = "  +3.60%  "
$ws.Range("D18").Value = "2.582.44"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").Value = "'338.67"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("D22").Value = "'6.26"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "'64.58"
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("D25").Value = "'0.457"
$ws.Range("E25").Value = "  +6.75%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("D29").Value = "0.0₃0781"
$ws.Range("E29").Value = "  +3.14%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'158.78"
$ws.Range("E33").Value = "  +2.99%  "
$ws.Range("D34").Value = "'19.02"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Value = "'4.03"
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("D37").Value = "'0.873"
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("D38").Value = "'0.874"
$ws.Range("E38").Value = "  -4.00%  "
$ws.Range("D39").Value = "'37.21"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'294.05"
$ws.Range("E41").Value = "  +4.15%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'3.68"
$ws.Range("E42").Value = "  +2.39%  "
$ws.Range("D44").Value = "'0.0976"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.594"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'128.06"
$ws.Range("E46").Value = "  +8.21%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "'19.20"
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").Value = "1.950.15"
$ws.Range("E51").Value = "  +0.01%  "
